$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-07 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-08 Saturday", 2) | Out-Null
$d.Content.Find.Execute("18÷6=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷3=4, 1", 2) | Out-Null
$d.Content.Find.Execute("56÷9=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "53÷7=7, 4", 2) | Out-Null
$d.Content.Find.Execute("12÷7=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "61÷7=8, 5", 2) | Out-Null
$d.Content.Find.Execute("13÷8=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=9, 6", 2) | Out-Null
$d.Content.Find.Execute("21÷3=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "63÷3=21, 0", 2) | Out-Null
$d.Content.Find.Execute("81÷3=27, 0", $true, $false, $false, $false, $false, $true, 1, $false, "14÷5=2, 4", 2) | Out-Null
$d.Content.Find.Execute("27÷2=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "30÷7=4, 2", 2) | Out-Null
$d.Content.Find.Execute("23÷3=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "70÷9=7, 7", 2) | Out-Null
$d.Content.Find.Execute("58÷3=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "62÷2=31, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷5=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=12, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷4=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=6, 4", 2) | Out-Null
$d.Content.Find.Execute("43÷7=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷3=13, 1", 2) | Out-Null
$d.Content.Find.Execute("12÷3=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "12÷5=2, 2", 2) | Out-Null
$d.Content.Find.Execute("65÷6=10, 5", $true, $false, $false, $false, $false, $true, 1, $false, "10÷8=1, 2", 2) | Out-Null
$d.Content.Find.Execute("34÷4=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "14÷8=1, 6", 2) | Out-Null
$d.Content.Find.Execute("36÷7=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "95÷2=47, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷4=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "89÷3=29, 2", 2) | Out-Null
$d.Content.Find.Execute("44÷5=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "29÷6=4, 5", 2) | Out-Null
$d.Content.Find.Execute("57÷4=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "59÷3=19, 2", 2) | Out-Null
$d.Content.Find.Execute("33÷8=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=18, 2", 2) | Out-Null
$d.Content.Find.Execute("79÷5=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "60÷5=12, 0", 2) | Out-Null
$d.Content.Find.Execute("75÷2=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷5=9, 1", 2) | Out-Null
$d.Content.Find.Execute("20÷9=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "82÷2=41, 0", 2) | Out-Null
$d.Content.Find.Execute("40÷9=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "79÷7=11, 2", 2) | Out-Null
$d.Content.Find.Execute("42÷2=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "16÷6=2, 4", 2) | Out-Null
